$d = $word.ActiveDocument

$pairs = @(
    @{old="65×88=5720"; new="76×28=2128"},
    @{old="59×33=1947"; new="82×62=5084"},
    @{old="68×11=748";  new="66×39=2574"},
    @{old="88×86=7568"; new="64×82=5248"},
    @{old="88×39=3432"; new="15×18=270"},
    @{old="64×84=5376"; new="17×38=646"},
    @{old="62×52=3224"; new="94×38=3572"},
    @{old="39×50=1950"; new="31×28=868"},
    @{old="84×49=4116"; new="16×77=1232"},
    @{old="32×87=2784"; new="73×19=1387"},
    @{old="32×51=1632"; new="40×28=1120"},
    @{old="70×50=3500"; new="40×68=2720"},
    @{old="52×55=2860"; new="58×60=3480"},
    @{old="25×42=1050"; new="32×52=1664"},
    @{old="47×27=1269"; new="36×14=504"},
    @{old="87×52=4524"; new="19×67=1273"},
    @{old="63×35=2205"; new="62×80=4960"},
    @{old="81×22=1782"; new="16×94=1504"},
    @{old="45×34=1530"; new="61×68=4148"},
    @{old="40×46=1840"; new="15×76=1140"},
    @{old="24×15=360";  new="39×98=3822"},
    @{old="43×47=2021"; new="63×29=1827"},
    @{old="17×93=1581"; new="81×92=7452"},
    @{old="72×65=4680"; new="92×23=2116"},
    @{old="31×83=2573"; new="62×55=3410"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
